# Append the new 2026-01-16 10:00:11 monitoring row (row 20) to the
# NIFTY_Options_Analysis sheet. The new row reuses the exact same
# per-column formatting (fill/font/border/alignment/number format)
# as the previous data row (row 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 19
$dstRow = 20
$srcRange = "A" + $srcRow + ":AE" + $srcRow
$dstRange = "A" + $dstRow + ":AE" + $dstRow

# 1) Clone row 19's formatting onto row 20 first, so every destination
#    cell already carries the right style before any value lands in it.
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial(-4122)   # xlPasteFormats

# 2) The new row's values, keyed by column letter.
$values = @{
    "A"  = "2026-01-16"
    "B"  = "10:00:11"
    "C"  = "AVOID"
    "D"  = "AVOID"
    "E"  = "100%"
    "F"  = "TRADEABLE"
    "G"  = 0
    "H"  = 25739.05
    "I"  = 11.16
    "J"  = -0.21
    "K"  = 0
    "L"  = 21.1
    "M"  = "UNKNOWN"
    "N"  = 0
    "O"  = "UNKNOWN"
    "P"  = 0
    "Q"  = 0
    "R"  = 0
    "S"  = 0
    "T"  = "NONE"
    "U"  = ""
    "V"  = 0
    "W"  = 0
    "X"  = 0
    "Y"  = 0
    "Z"  = 0
    "AA" = 0
    "AB" = 0
    "AC" = "HARD VETO: CPR TRENDING DAY: Price 25739.05 above TC 25676.35 - BULLISH TRENDING DAY likely"
    "AD" = "CPR TRENDING DAY: Price 25739.05 above TC 25676.35 - BULLISH TRENDING DAY likely"
    "AE" = "Yes"
}

# Columns whose literal text reads like a date ("2026-01-16"), a
# percentage ("100%"), or is simply empty ("") need a leading
# quote-prefix when assigned - otherwise the value setter either
# auto-converts the text into a date serial / percentage number, or
# silently drops an empty-string assignment instead of leaving a
# blank text cell behind. The quote-prefix keeps the literal text
# intact (as a genuine text cell) but bumps the cell onto a brand
# new ad-hoc style, so row 19's formatting is stamped on again
# afterwards to bring every cell in the row back onto the shared
# per-column styles.
$textForced = @("A", "E", "U")

foreach ($col in $values.Keys) {
    $cell = $ws.Range($col + $dstRow)
    if ($textForced -contains $col) {
        $cell.Value = "'" + $values[$col]
    } else {
        $cell.Value = $values[$col]
    }
}

# 3) Re-apply row 19's formatting so every cell - including the ones
#    touched by the quote-prefix workaround above - ends up back on
#    the correct shared style.
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial(-4122)   # xlPasteFormats
